$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values for rows 2-11
$ws.Range("C2").Value = -0.2112792329596069
$ws.Range("D2").Value = 0.8346142853460654

$ws.Range("C3").Value = 0.8749341390389186
$ws.Range("D3").Value = 0.3910618507920902

$ws.Range("C4").Value = -0.02303693466243946
$ws.Range("D4").Value = 0.9818284854960273

$ws.Range("C5").Value = 0.03666118379617099
$ws.Range("D5").Value = 0.9710857866805345

$ws.Range("C6").Value = 0.9345532631462787
$ws.Range("D6").Value = 0.3601637795523547

$ws.Range("C7").Value = 0.2270016199225518
$ws.Range("D7").Value = 0.8225185970915152

$ws.Range("C8").Value = 0.2119308746452045
$ws.Range("D8").Value = 0.8341121038847605

$ws.Range("C9").Value = -0.933357001044356
$ws.Range("D9").Value = 0.3607672231615782

$ws.Range("C10").Value = -0.65224488974074
$ws.Range("D10").Value = 0.5210004191458175

$ws.Range("C11").Value = 0.07578543954818262
$ws.Range("D11").Value = 0.9402747530856583
